$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 48, columns C..I) that was added under
# the "4" (April) month group.
$ws.Range("C48").Value = 768
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 22
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 300
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0

# Extend the shared "diff" formula down to the new row.
$ws.Range("J48").Formula = "=(H48+I48)-(C48+D48+E48+F48+G48)"

# Move the active cell / selection to the new last row, matching the
# author's saved cursor position.
$ws.Range("J48").Select()
